$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaStart = $metaRange.Start

$metaRange.InsertAfter("Meta description: Find out why you might want to play 50 Lions for free, including unique bonus features and high-quality graphics.")

# Make just the "Meta description" label bold.
$metaBoldRange = $d.Range($metaStart, $metaStart + 16)
$metaBoldRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated title paragraph that was sitting right
#    before the final (italic) paragraph near the end of the document.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the final paragraph's text (the italic meta-description
#    placeholder) with the new AI image-generation prompt.
# ------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count2)
$lastRange = $lastPara.Range
# Exclude the trailing paragraph mark so the assignment replaces the
# run's text in place instead of inserting before it.
$lastTextRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTextRange.Text = 'Create a feature image for the game "50 Lions" that embodies the excitement and adventure of the African savannah. The image should be in cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be surrounded by iconic African animals, including lions, zebras, giraffes, and elephants. The background should be a beautiful sunset over the savannah landscape with an Acacia tree on one side. The image should capture the essence of the game, which is to embark on a virtual safari adventure and win big with the help of bonus features and free spins.'
